$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the data for day 9 (row 17), which was previously blank
$ws.Range("A17").Value = 9
$ws.Range("D17").Value = "Finalizare 5.1 si 5.2"
$ws.Range("B17").Value = "18/6/2024"
$ws.Range("C17").Value = 3.3

# Insert 4 new blank rows before the totals section (old rows 19-21 shift to 23-25)
$ws.Rows("19:22").Insert()

# Update the total-hours formula to cover the newly inserted rows
$ws.Range("D23").Formula = "=SUM(C3:C22)"
$ws.Rows(23).AutoFit()

# Update the selected cell to match the saved view state
$ws.Range("C19").Select()

Write-Host "done"
